# The workbook is already open; grab references.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cell D1 was edited from " Humidity" (leading space) to "Humidity".
$ws.Range("D1").Value = "Humidity"

# Columns C and D were given explicit custom widths.
$ws.Columns("C").ColumnWidth = 14
$ws.Columns("D").ColumnWidth = 20.8

# The active selection on the sheet moved to G14.
$ws.Range("G14").Select() | Out-Null
